# Apply retrained model numbers to the yolov3_summary worksheet.
# Addresses issue #15 - update output channel (O, and L for post-processing
# rows) counts from 30 to 21, and recompute the per-row timings (S column)
# and grand total (S114) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 89, 98, 107: output conv layers - only column O changes (30 -> 21)
# plus their time[us] values in column S.
$ws.Cells.Item(89, 15).Value  = 21    # O89
$ws.Cells.Item(89, 19).Value  = 125   # S89

$ws.Cells.Item(98, 15).Value  = 21    # O98
$ws.Cells.Item(98, 19).Value  = 226   # S98

$ws.Cells.Item(107, 15).Value = 21    # O107
$ws.Cells.Item(107, 19).Value = 469   # S107

# Rows 108-113: post-processing transpose/cast rows - both L (in channels)
# and O (out channels) change from 30 -> 21, plus their time[us] values.
$ws.Cells.Item(108, 12).Value = 21    # L108
$ws.Cells.Item(108, 15).Value = 21    # O108
$ws.Cells.Item(108, 19).Value = 4     # S108

$ws.Cells.Item(109, 12).Value = 21    # L109
$ws.Cells.Item(109, 15).Value = 21    # O109
$ws.Cells.Item(109, 19).Value = 3     # S109

$ws.Cells.Item(110, 12).Value = 21    # L110
$ws.Cells.Item(110, 15).Value = 21    # O110
$ws.Cells.Item(110, 19).Value = 14    # S110

$ws.Cells.Item(111, 12).Value = 21    # L111
$ws.Cells.Item(111, 15).Value = 21    # O111
$ws.Cells.Item(111, 19).Value = 12    # S111

$ws.Cells.Item(112, 12).Value = 21    # L112
$ws.Cells.Item(112, 15).Value = 21    # O112
$ws.Cells.Item(112, 19).Value = 54    # S112

$ws.Cells.Item(113, 12).Value = 21    # L113
$ws.Cells.Item(113, 15).Value = 21    # O113
$ws.Cells.Item(113, 19).Value = 45    # S113

# Row 114: grand total time[us]
$ws.Cells.Item(114, 19).Value = 189131   # S114
